# Add a progress-update blurb (plus a blank spacer line) after the final
# screenshot image at the end of the report, matching the author's
# "Wall and ceiling collision added" commit.

$d = $word.ActiveDocument

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Add-CleanParagraph([string]$innerXml) {
    # Appends a brand-new paragraph at the very end of the document body.
    # A plain InsertParagraphAfter() leaves a stray empty <w:r/> behind in
    # this engine, so instead we grow the document with a scratch trailing
    # paragraph and then overwrite it via InsertXML with exactly the OOXML
    # we want (no run at all for the blank spacer paragraph, a single
    # <w:r><w:t>...</w:t></w:r> for the text ones) - matching how Word
    # itself would have serialized a freshly typed paragraph.
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastRange = $lastPara.Range
    $lastRange.Collapse(0)            # wdCollapseEnd
    $lastRange.InsertParagraphAfter()

    $scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $scratchRange = $scratchPara.Range
    $scratchRange.Collapse(1)         # wdCollapseStart
    $scratchRange.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

$tabsPPr = "<w:pPr><w:tabs><w:tab w:val='left' w:pos='5629'/></w:tabs></w:pPr>"

# 1. blank spacer paragraph (tab stop only, no run)
Add-CleanParagraph($tabsPPr)

# 2. "A lot of progress..." paragraph
Add-CleanParagraph($tabsPPr + "<w:r><w:t>A lot of progress, and now the player can travel along full loops.</w:t></w:r>")

# 3. "Next up is wall/ceiling collision..." paragraph
Add-CleanParagraph($tabsPPr + "<w:r><w:t>Next up is wall/ceiling collision. Then movement.</w:t></w:r>")
